$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 542.6842
$ws.Range("J17").Value = 287.4
$ws.Range("L17").Value = 862.1999999999999
$ws.Range("N17").Value = -1198.2
$ws.Range("H28").Value = 925.4583
$ws.Range("I28").Value = 917.75
$ws.Range("J28").Value = 964
$ws.Range("K28").Value = 917.75
$ws.Range("L28").Value = 964
$ws.Range("M28").Value = -432.75
$ws.Range("N28").Value = -1934
$ws.Range("H76").Value = 3101.5107
$ws.Range("I76").Value = 2994.6191
$ws.Range("K76").Value = 2994.6191
$ws.Range("M76").Value = -2679.6191
$ws.Range("H79").Value = 3101.5107
$ws.Range("I79").Value = 2994.6191
$ws.Range("K79").Value = 2994.6191
$ws.Range("M79").Value = -1902.6191
$ws.Range("H107").Value = 1077.48
$ws.Range("I107").Value = 1101.5416
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 1101.5416
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 818.4584
$ws.Range("N107").Value = -4340
$ws.Range("H111").Value = 1333.3334
$ws.Range("I111").Value = 750
$ws.Range("J111").Value = 2500
$ws.Range("K111").Value = 2250
$ws.Range("L111").Value = 7500
$ws.Range("M111").Value = 817
$ws.Range("N111").Value = -13634
$ws.Range("H112").Value = 1329.3469
$ws.Range("J112").Value = 1346.625
$ws.Range("L112").Value = 4039.875
$ws.Range("N112").Value = -6255.875
$ws.Range("H137").Value = 746962.9
$ws.Range("I137").Value = 1908024.6
$ws.Range("J137").Value = 2692.4614
$ws.Range("K137").Value = 5724073.800000001
$ws.Range("L137").Value = 8077.3842
$ws.Range("M137").Value = -5721523.800000001
$ws.Range("N137").Value = -13177.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5138.603
$ws.Range("I32").Value = 5573.7144
$ws.Range("J32").Value = 4268.381
$ws.Range("K32").Value = 5573.7144
$ws.Range("L32").Value = 4268.381
$ws.Range("M32").Value = -5286.7144
$ws.Range("N32").Value = -4842.381
$ws.Range("H45").Value = 4418.4287
$ws.Range("I45").Value = 3232.5
$ws.Range("J45").Value = 5999.6665
$ws.Range("K45").Value = 3232.5
$ws.Range("L45").Value = 5999.6665
$ws.Range("M45").Value = -2855.5
$ws.Range("N45").Value = -6753.6665
$ws.Range("H122").Value = 4256
$ws.Range("I122").Value = 1800.3334
$ws.Range("J122").Value = 7202.8
$ws.Range("K122").Value = 5401.0002
$ws.Range("L122").Value = 21608.4
$ws.Range("M122").Value = -2951.0002
$ws.Range("N122").Value = -26508.4
$ws.Range("H132").Value = 2519.0322
$ws.Range("I132").Value = 1278.5264
$ws.Range("J132").Value = 4483.1665
$ws.Range("K132").Value = 3835.5792
$ws.Range("L132").Value = 13449.4995
$ws.Range("M132").Value = -1305.5792
$ws.Range("N132").Value = -18509.4995
$ws.Range("H139").Value = 43076.25
$ws.Range("J139").Value = 43076.25
$ws.Range("L139").Value = 43076.25
$ws.Range("N139").Value = -53356.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 326.9375
$ws.Range("J80").Value = 262
$ws.Range("L80").Value = 262
$ws.Range("N80").Value = -2258
$ws.Range("H83").Value = 326.9375
$ws.Range("J83").Value = 262
$ws.Range("L83").Value = 1310
$ws.Range("N83").Value = -11294
$ws.Range("H99").Value = 1407.7778
$ws.Range("I99").Value = 1281.4286
$ws.Range("J99").Value = 1850
$ws.Range("K99").Value = 1281.4286
$ws.Range("L99").Value = 1850
$ws.Range("M99").Value = 216.5714
$ws.Range("N99").Value = -4846
$ws.Range("H134").Value = 2400.7856
$ws.Range("I134").Value = 1307.5358
$ws.Range("J134").Value = 4587.2856
$ws.Range("K134").Value = 3922.6074
$ws.Range("L134").Value = 13761.8568
$ws.Range("M134").Value = -1387.6074
$ws.Range("N134").Value = -18831.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 170721.28
$ws.Range("I31").Value = 386773.7
$ws.Range("J31").Value = 2680.5112
$ws.Range("K31").Value = 386773.7
$ws.Range("L31").Value = 2680.5112
$ws.Range("M31").Value = -386478.7
$ws.Range("N31").Value = -3270.5112
$ws.Range("H34").Value = 170721.28
$ws.Range("I34").Value = 386773.7
$ws.Range("J34").Value = 2680.5112
$ws.Range("K34").Value = 386773.7
$ws.Range("L34").Value = 2680.5112
$ws.Range("M34").Value = -386571.7
$ws.Range("N34").Value = -3084.5112
$ws.Range("H58").Value = 2079.9778
$ws.Range("I58").Value = 1227.7667
$ws.Range("J58").Value = 3784.4
$ws.Range("K58").Value = 1227.7667
$ws.Range("L58").Value = 3784.4
$ws.Range("M58").Value = -1024.7667
$ws.Range("N58").Value = -4190.4
$ws.Range("H122").Value = 2778.5
$ws.Range("I122").Value = 2077.4285
$ws.Range("J122").Value = 3760
$ws.Range("K122").Value = 6232.2855
$ws.Range("L122").Value = 11280
$ws.Range("M122").Value = -3782.2855
$ws.Range("N122").Value = -16180
$ws.Range("H132").Value = 3297.6155
$ws.Range("I132").Value = 2553.913
$ws.Range("K132").Value = 7661.739
$ws.Range("M132").Value = -5131.739
$ws.Range("H134").Value = 1544.4062
$ws.Range("I134").Value = 1003.13043
$ws.Range("J134").Value = 2927.6667
$ws.Range("K134").Value = 3009.39129
$ws.Range("L134").Value = 8783.000100000001
$ws.Range("M134").Value = -474.39129
$ws.Range("N134").Value = -13853.0001
$ws.Range("H136").Value = 2079.9778
$ws.Range("I136").Value = 1227.7667
$ws.Range("J136").Value = 3784.4
$ws.Range("K136").Value = 3683.300099999999
$ws.Range("L136").Value = 11353.2
$ws.Range("M136").Value = -1133.300099999999
$ws.Range("N136").Value = -16453.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2693.75
$ws.Range("I5").Value = 1620.2
$ws.Range("K5").Value = 4860.6
$ws.Range("M5").Value = -4748.6
$ws.Range("H92").Value = 1315.8
$ws.Range("I92").Value = 474.5
$ws.Range("J92").Value = 1876.6666
$ws.Range("K92").Value = 1423.5
$ws.Range("L92").Value = 5629.9998
$ws.Range("M92").Value = -175.5
$ws.Range("N92").Value = -8125.9998
$ws.Range("H98").Value = 563.8333
$ws.Range("J98").Value = 627.5714
$ws.Range("L98").Value = 1882.7142
$ws.Range("N98").Value = -4878.7142
$ws.Range("H113").Value = 464.75
$ws.Range("I113").Value = 473.96155
$ws.Range("J113").Value = 451.44446
$ws.Range("K113").Value = 1421.88465
$ws.Range("L113").Value = 1354.33338
$ws.Range("M113").Value = 748.11535
$ws.Range("N113").Value = -5694.33338
$ws.Range("H121").Value = 1636.4557
$ws.Range("I121").Value = 515.5
$ws.Range("J121").Value = 1877.8923
$ws.Range("K121").Value = 1546.5
$ws.Range("L121").Value = 5633.6769
$ws.Range("M121").Value = -236.5
$ws.Range("N121").Value = -8253.6769
$ws.Range("H131").Value = 820.59
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 827.1326
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2481.3978
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12561.3978
$ws.Range("H132").Value = 2377.125
$ws.Range("J132").Value = 2749.6924
$ws.Range("L132").Value = 24747.2316
$ws.Range("N132").Value = -29807.2316
$ws.Range("H135").Value = 2693.75
$ws.Range("I135").Value = 1620.2
$ws.Range("K135").Value = 14581.8
$ws.Range("M135").Value = -12046.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 29490
$ws.Range("J112").Value = 29490
$ws.Range("L112").Value = 29490
$ws.Range("N112").Value = -31706
$ws.Range("H122").Value = 3403.8462
$ws.Range("I122").Value = 1990
$ws.Range("K122").Value = 5970
$ws.Range("M122").Value = -3520
$ws.Range("H132").Value = 4726.3213
$ws.Range("I132").Value = 3732.2144
$ws.Range("J132").Value = 5720.4287
$ws.Range("K132").Value = 11196.6432
$ws.Range("L132").Value = 17161.2861
$ws.Range("M132").Value = -8666.643199999999
$ws.Range("N132").Value = -22221.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4577.8887
$ws.Range("I40").Value = 3644.375
$ws.Range("J40").Value = 5935.727
$ws.Range("K40").Value = 3644.375
$ws.Range("L40").Value = 5935.727
$ws.Range("M40").Value = -3508.375
$ws.Range("N40").Value = -6207.727
$ws.Range("H100").Value = 4436.091
$ws.Range("I100").Value = 1759.4
$ws.Range("J100").Value = 6666.6665
$ws.Range("K100").Value = 1759.4
$ws.Range("L100").Value = 6666.6665
$ws.Range("M100").Value = -1218.4
$ws.Range("N100").Value = -7748.6665
$ws.Range("H122").Value = 4868.1333
$ws.Range("I122").Value = 4273.143
$ws.Range("J122").Value = 5388.75
$ws.Range("K122").Value = 12819.429
$ws.Range("L122").Value = 16166.25
$ws.Range("M122").Value = -10369.429
$ws.Range("N122").Value = -21066.25
$ws.Range("H132").Value = 4352.227
$ws.Range("I132").Value = 3123.3333
$ws.Range("K132").Value = 9369.999899999999
$ws.Range("M132").Value = -6839.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9264785
$ws.Range("I132").Value = 8764.214
$ws.Range("J132").Value = 15154980
$ws.Range("K132").Value = 26292.642
$ws.Range("L132").Value = 45464940
$ws.Range("M132").Value = -23762.642
$ws.Range("N132").Value = -45470000
$ws.Range("H136").Value = 6484.793
$ws.Range("I136").Value = 5596.05
$ws.Range("J136").Value = 8459.777
$ws.Range("K136").Value = 16788.15
$ws.Range("L136").Value = 25379.331
$ws.Range("M136").Value = -14238.15
$ws.Range("N136").Value = -30479.331
